$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" property to the new generation time
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# ------------------------------------------------------------------
# 2) Elements sheet: append a new element row describing the
#    ActiviteSociale.EntiteGeographique reference element
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Duplicate the last existing data row (20) into the new row (21) so the
# new row inherits the same cell style/shape (borders, wrap, blanks) as
# every other element row, then overwrite the cells that actually differ.
$ws.Range("A20:AJ20").Copy($ws.Range("A21:AJ21"))

$path = "ActiviteSociale.EntiteGeographique"

$ws.Cells.Item(21, 1).Value = $path   # ID
$ws.Cells.Item(21, 2).Value = $path   # Path

# Min / Max (column F / G) and Base Min / Base Max (AG / AH) are stored as
# text "1" (not a number) in this sheet, same as every other row - copy a
# cell that already holds a text "1" so the type matches instead of
# letting a plain Value assignment coerce it to a numeric cell.
$ws.Range("G3").Copy($ws.Range("F21"))
$ws.Range("G3").Copy($ws.Range("G21"))
$ws.Range("G3").Copy($ws.Range("AG21"))
$ws.Range("G3").Copy($ws.Range("AH21"))

# Type(s) / Short / Definition
$ws.Cells.Item(21, 11).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteGeographique`n"
$ws.Cells.Item(21, 12).Value = "Lien vers la classe EntiteGeographique"
$ws.Cells.Item(21, 13).Value = "Lien vers la classe EntiteGeographique"

# Base Path
$ws.Cells.Item(21, 32).Value = $path

# The copied row 20 carried Binding Strength / Binding Value Set
# ("preferred" / a TRE url) in columns X and Z - the new row has neither,
# so clear those back out.
$ws.Range("X21").ClearContents()
$ws.Range("Z21").ClearContents()
